$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace all values of 2 with 0 in columns L, M and N (rows 2-112),
# i.e. recode the "2" responses of pubfund1 / programpart / startupcenter to 0.
$rng = $ws.Range("L2:N112")
$rng.Replace(2, 0, [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)

# Update the active view/selection on the sheet
$ws.Activate()
$ws.Range("L1:N1048576").Select()
$excel.ActiveWindow.ScrollRow = 22

$wb.Save()
